$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Prefix the date (with a line break) onto the header cells D1:J1
$ws.Range("D1").Value = "2月16日`n马太福音13-14"
$ws.Range("E1").Value = "2月17日`n马太福音15-16"
$ws.Range("F1").Value = "2月18日`n马太福音17-18"
$ws.Range("G1").Value = "2月19日`n马太福音19-20"
$ws.Range("H1").Value = "2月20日`n马太福音21-22"
$ws.Range("I1").Value = "2月21日`n马太福音23-24"
$ws.Range("J1").Value = "2月22日`n马太福音25-26"

# 2) Make row 1 taller (30pt) and mark it as a custom height
$ws.Range("A1:J1").EntireRow.RowHeight = 30

# 3) Turn on word-wrap for all the cells in the sheet (this affects the two
#    shared cell styles used throughout the sheet, matching the xf wrapText
#    change in the stylesheet)
$ws.Cells.WrapText = $true

# 4) Shrink the first conditional formatting rule's range from B2:C100 to B2:B100
$fc = $ws.Range("B2").FormatConditions.Item(1)
$fc.ModifyAppliesToRange($ws.Range("B2:B100"))
